$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "49.183.13"
$ws.Cells.Item(2, 5).Value = "  -1.39%  "

Set-TextValue 3 4 "2.633.17"
$ws.Cells.Item(3, 5).Value = "  +0.26%  "

Set-TextValue 4 4 "0.999"
$ws.Cells.Item(4, 5).Value = "  +0.10%  "

Set-TextValue 5 4 "111.60"
$ws.Cells.Item(5, 5).Value = "  +1.09%  "

Set-TextValue 6 4 "323.16"

Set-TextValue 7 4 "0.526"
$ws.Cells.Item(7, 5).Value = "  -2.00%  "

Set-TextValue 8 4 "0.999"
$ws.Cells.Item(8, 5).Value = "  +0.03%  "

Set-TextValue 9 4 "0.544"
$ws.Cells.Item(9, 5).Value = "  -3.68%  "

Set-TextValue 10 4 "39.81"
$ws.Cells.Item(10, 5).Value = "  -3.04%  "

Set-TextValue 11 4 "19.81"
$ws.Cells.Item(11, 5).Value = "  -4.19%  "

Set-TextValue 12 4 "0.0811"
$ws.Cells.Item(12, 5).Value = "  -1.96%  "

$ws.Cells.Item(13, 5).Value = "  +0.19%  "

Set-TextValue 14 4 "7.26"
$ws.Cells.Item(14, 5).Value = "  -1.22%  "

Set-TextValue 15 4 "3.034.65"
$ws.Cells.Item(15, 5).Value = "  +0.33%  "

Set-TextValue 16 4 "2.619.28"
$ws.Cells.Item(16, 5).Value = "  +0.55%  "

Set-TextValue 17 4 "0.862"
$ws.Cells.Item(17, 5).Value = "  -1.49%  "

Set-TextValue 18 4 "49.126.21"
$ws.Cells.Item(18, 5).Value = "  -1.39%  "

Set-TextValue 19 4 "2.99"
$ws.Cells.Item(19, 5).Value = "  -3.62%  "

Set-TextValue 20 4 "12.89"
$ws.Cells.Item(20, 5).Value = "  -3.77%  "

Set-TextValue 21 4 "6.70"
$ws.Cells.Item(21, 5).Value = "  -1.59%  "

Set-TextValue 22 4 "0.0₃0947"
$ws.Cells.Item(22, 5).Value = "  -1.13%  "

Set-TextValue 23 4 "269.89"
$ws.Cells.Item(23, 5).Value = "  -4.52%  "

Set-TextValue 24 4 "68.55"
$ws.Cells.Item(24, 5).Value = "  -5.99%  "

Set-TextValue 25 4 "2.55"
$ws.Cells.Item(25, 5).Value = "  -1.19%  "

Set-TextValue 26 4 "26.15"
$ws.Cells.Item(26, 5).Value = "  -2.40%  "

$ws.Cells.Item(27, 5).Value = "  +0.05%  "

Set-TextValue 28 4 "10.06"
$ws.Cells.Item(28, 5).Value = "  +0.60%  "

$ws.Cells.Item(29, 5).Value = "  -0.82%  "

Set-TextValue 30 4 "35.06"
$ws.Cells.Item(30, 5).Value = "  -3.37%  "

Set-TextValue 31 4 "0.138"
$ws.Cells.Item(31, 5).Value = "  -4.86%  "

Set-TextValue 32 4 "49.48"
$ws.Cells.Item(32, 5).Value = "  -0.28%  "

Set-TextValue 33 4 "5.50"
$ws.Cells.Item(33, 5).Value = "  +0.16%  "

$ws.Cells.Item(34, 5).Value = "  -0.21%  "

Set-TextValue 37 4 "5.00"
$ws.Cells.Item(37, 5).Value = "  +4.78%  "

$ws.Cells.Item(38, 5).Value = "  -1.09%  "

Set-TextValue 39 4 "3.14"
$ws.Cells.Item(39, 5).Value = "  +1.79%  "

Set-TextValue 40 4 "127.01"
$ws.Cells.Item(40, 5).Value = "  +2.56%  "

$ws.Cells.Item(41, 5).Value = "  -1.71%  "

Set-TextValue 42 4 "22.27"
$ws.Cells.Item(42, 5).Value = "  -3.27%  "

$ws.Cells.Item(43, 5).Value = "  -4.36%  "

$ws.Cells.Item(44, 5).Value = "  -0.16%  "

Set-TextValue 45 4 "2.066.66"
$ws.Cells.Item(45, 5).Value = "  +0.55%  "

Set-TextValue 46 4 "2.17"
$ws.Cells.Item(46, 5).Value = "  +6.51%  "

Set-TextValue 47 4 "3.26"
$ws.Cells.Item(47, 5).Value = "  -2.88%  "

$ws.Cells.Item(48, 5).Value = "  -2.46%  "

Set-TextValue 49 4 "8.86"
$ws.Cells.Item(49, 5).Value = "  -2.56%  "

$ws.Cells.Item(35, 2).Value = "Celestia"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue 35 4 "19.06"
$ws.Cells.Item(35, 5).Value = "  -3.42%  "

$ws.Cells.Item(36, 2).Value = "Hedera"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue 36 4 "0.0797"
$ws.Cells.Item(36, 5).Value = "  -0.01%  "

$ws.Cells.Item(50, 2).Value = "THORChain"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue 50 4 "5.19"
$ws.Cells.Item(50, 5).Value = "  -3.76%  "

$ws.Cells.Item(51, 2).Value = "MultiversX"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue 51 4 "58.71"
$ws.Cells.Item(51, 5).Value = "  +1.46%  "
